# Commit: "template and stack mods pulling back from 8.1 configs"
#
# Changes applied to the "values" sheet:
#   1. PANORAMA_TYPE value: "cloud" -> "static"
#   2. PANORAMA_IP value:   "192.168.55.7" -> "192.168.55.8"
#   3. A new row is inserted right after CONFIG_EXPORT_IP (before STACK) that
#      defines a new TEMPLATE variable for Panorama:
#         TEMPLATE | sample_template | Template name for Panorama
#   4. MGMT_TYPE value: "dhcp-cloud" -> "static", and its description is
#      updated from "firewall management IP type (static or dhcp-cloud)" to
#      "firewall management IP type (static or dhcp-client)"
#
# All of the "set commands" sheet formulas that reference 'values'!B<n> shift
# automatically because row 8 gets physically inserted into the "values"
# sheet (Excel keeps formula references in sync across sheets on insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# 1. PANORAMA_TYPE (row 3) value cloud -> static
$ws.Range("B3").Value = "static"

# 2. PANORAMA_IP (row 4) value 192.168.55.7 -> 192.168.55.8
$ws.Range("B4").Value = "192.168.55.8"

# 3. Insert a new row before the current STACK row (row 8) and populate it
#    with the new TEMPLATE variable. This naturally shifts STACK and every
#    row below it down by one, and Excel updates every cross-sheet formula
#    reference to 'values'!B<n> accordingly.
$ws.Range("A8").EntireRow.Insert()
$ws.Range("A8").Value = "TEMPLATE"
$ws.Range("B8").Value = "sample_template"
$ws.Range("C8").Value = "Template name for Panorama"

# 4. MGMT_TYPE row (now row 12 after the insert above) value dhcp-cloud -> static
#    and update its description text.
$ws.Range("B12").Value = "static"
$ws.Range("C12").Value = "firewall management IP type (static or dhcp-client)"
